# Reprocess D2893, Edi06 and Edi09
# Update detection-limit values in both worksheets of the Edi09 calczaf_outputs workbook.

$wb = $excel.ActiveWorkbook

# Column order used across the data rows (B..R)
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")

# ---- Sheet 1: Edi09_2_bg_detlim ----
$ws1 = $wb.Worksheets.Item("Edi09_2_bg_detlim")

$row2vals = @(0.023,0.021,0.025,0.026,0.025,0.02,0.028,0.026,0.027,0.024,0.03,0.028,0.027,0.025,0.003,0.02,0.03)
$row8vals = @(0.032,0.03,0.035,0.036,0.035,0.029,0.039,0.036,0.038,0.034,0.042,0.04,0.038,0.036,0.004,0.029,0.042)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws1.Range($cols[$i] + "2").Value = $row2vals[$i]
    $ws1.Range($cols[$i] + "5").Value = $row2vals[$i]
    $ws1.Range($cols[$i] + "8").Value = $row8vals[$i]
}

# ---- Sheet 2: Edi09_3_bg_apf_detlim ----
$ws2 = $wb.Worksheets.Item("Edi09_3_bg_apf_detlim")

$row2vals2 = @(0.027,0.025,0.03,0.031,0.03,0.025,0.033,0.031,0.032,0.029,0.036,0.034,0.033,0.03,0.003,0.025,0.036)
$row8vals2 = @(0.038,0.036,0.042,0.044,0.042,0.035,0.047,0.044,0.045,0.04,0.05,0.048,0.046,0.043,0.005,0.035,0.05)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws2.Range($cols[$i] + "2").Value = $row2vals2[$i]
    $ws2.Range($cols[$i] + "5").Value = $row2vals2[$i]
    $ws2.Range($cols[$i] + "8").Value = $row8vals2[$i]
}
